$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1904761904761905
$ws.Range("C2").Value = 0.575091575091575
$ws.Range("J2").Value = 0.01831501831501832
$ws.Range("P2").Value = 0.1172161172161172
$ws.Range("S2").Value = 0.0989010989010989
$ws.Range("C3").Value = 0.03048780487804878
$ws.Range("J3").Value = 0.0426829268292683
$ws.Range("P3").Value = 0.7439024390243902
$ws.Range("S3").Value = 0.1829268292682927
$ws.Range("J4").Value = 0.04761904761904762
$ws.Range("P4").Value = 0.7142857142857143
$ws.Range("S4").Value = 0.2380952380952381
$ws.Range("B6").Value = 0.04273504273504274
$ws.Range("D6").Value = 0.02136752136752137
$ws.Range("F6").Value = 0.07264957264957266
$ws.Range("J6").Value = 0.2564102564102564
$ws.Range("O6").Value = 0.0170940170940171
$ws.Range("Q6").Value = 0.1709401709401709
$ws.Range("R6").Value = 0.05982905982905983
$ws.Range("S6").Value = 0.358974358974359
$ws.Range("B7").Value = 0.09047619047619047
$ws.Range("D7").Value = 0.02380952380952381
$ws.Range("F7").Value = 0.0761904761904762
$ws.Range("J7").Value = 0.1047619047619048
$ws.Range("O7").Value = 0.009523809523809525
$ws.Range("Q7").Value = 0.1714285714285714
$ws.Range("R7").Value = 0.09047619047619047
$ws.Range("S7").Value = 0.4333333333333333
$ws.Range("B8").Value = 0.08481262327416174
$ws.Range("D8").Value = 0.007889546351084813
$ws.Range("E8").Value = 0.001972386587771203
$ws.Range("F8").Value = 0.05522682445759369
$ws.Range("J8").Value = 0.1104536489151874
$ws.Range("O8").Value = 0.01577909270216963
$ws.Range("Q8").Value = 0.2169625246548323
$ws.Range("R8").Value = 0.1143984220907298
$ws.Range("S8").Value = 0.3925049309664694
$ws.Range("B9").Value = 0.07662835249042145
$ws.Range("D9").Value = 0.01532567049808429
$ws.Range("F9").Value = 0.08045977011494253
$ws.Range("J9").Value = 0.09961685823754789
$ws.Range("O9").Value = 0.01149425287356322
$ws.Range("Q9").Value = 0.1762452107279693
$ws.Range("R9").Value = 0.1187739463601533
$ws.Range("S9").Value = 0.421455938697318
$ws.Range("B10").Value = 0.09300648882480173
$ws.Range("D10").Value = 0.01802451333813987
$ws.Range("E10").Value = 0.002883922134102379
$ws.Range("F10").Value = 0.05912040374909877
$ws.Range("J10").Value = 0.1023792357606345
$ws.Range("O10").Value = 0.007930785868781542
$ws.Range("Q10").Value = 0.2278298485940879
$ws.Range("R10").Value = 0.08651766402307137
$ws.Range("S10").Value = 0.4023071377072819
$ws.Range("G11").Value = 0.1596091205211726
$ws.Range("J11").Value = 0.09771986970684039
$ws.Range("K11").Value = 0.2117263843648209
$ws.Range("L11").Value = 0.5179153094462541
$ws.Range("S11").Value = 0.01302931596091205
$ws.Range("G12").Value = 0.7682926829268293
$ws.Range("J12").Value = 0.1707317073170732
$ws.Range("K12").Value = 0.01219512195121951
$ws.Range("L12").Value = 0.02439024390243903
$ws.Range("S12").Value = 0.02439024390243903
$ws.Range("G13").Value = 0.875
$ws.Range("J13").Value = 0.125
$ws.Range("F15").Value = 0.008695652173913044
$ws.Range("H15").Value = 0.1782608695652174
$ws.Range("I15").Value = 0.1
$ws.Range("J15").Value = 0.4
$ws.Range("K15").Value = 0.03478260869565217
$ws.Range("M15").Value = 0.0391304347826087
$ws.Range("O15").Value = 0.05217391304347826
$ws.Range("S15").Value = 0.1869565217391304
$ws.Range("F16").Value = 0.03314917127071823
$ws.Range("H16").Value = 0.1933701657458564
$ws.Range("I16").Value = 0.1104972375690608
$ws.Range("J16").Value = 0.4143646408839779
$ws.Range("K16").Value = 0.09944751381215469
$ws.Range("M16").Value = 0.02209944751381215
$ws.Range("N16").Value = 0.005524861878453038
$ws.Range("O16").Value = 0.04972375690607735
$ws.Range("S16").Value = 0.0718232044198895
$ws.Range("F17").Value = 0.01291512915129151
$ws.Range("H17").Value = 0.1845018450184502
$ws.Range("I17").Value = 0.1107011070110701
$ws.Range("J17").Value = 0.4354243542435424
$ws.Range("K17").Value = 0.08118081180811808
$ws.Range("M17").Value = 0.01291512915129151
$ws.Range("N17").Value = 0.001845018450184502
$ws.Range("O17").Value = 0.07380073800738007
$ws.Range("S17").Value = 0.08671586715867159
$ws.Range("F18").Value = 0.008298755186721992
$ws.Range("H18").Value = 0.1493775933609958
$ws.Range("I18").Value = 0.1120331950207469
$ws.Range("J18").Value = 0.4356846473029046
$ws.Range("K18").Value = 0.1120331950207469
$ws.Range("M18").Value = 0.02489626556016597
$ws.Range("O18").Value = 0.07053941908713693
$ws.Range("S18").Value = 0.08713692946058091
$ws.Range("F19").Value = 0.0171003717472119
$ws.Range("H19").Value = 0.2223048327137546
$ws.Range("I19").Value = 0.1003717472118959
$ws.Range("J19").Value = 0.3836431226765799
$ws.Range("K19").Value = 0.1026022304832714
$ws.Range("M19").Value = 0.01933085501858736
$ws.Range("N19").Value = 0.0007434944237918215
$ws.Range("O19").Value = 0.07137546468401487
$ws.Range("S19").Value = 0.08252788104089219
